$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unneeded rows (old rows 7-22) by deleting entire rows,
# which shifts everything up and shrinks the used range.
$ws.Range("A7:A22").EntireRow.Delete()

# Consolidate each card's fields into a single Python-tuple-like string per row.
$ws.Range("A2").Value = "('Brave the Elements', ['{W}', 'Instant', 'Choose a color. White creatures you control gain protection from the chosen color until end of turn.'])"
$ws.Range("A3").Value = "('Day of Judgment', ['{2}{W}{W}', 'Sorcery', 'Destroy all creatures.'])"
$ws.Range("A4").Value = "('Doom Blade', ['{1}{B}', 'Instant', 'Destroy target nonblack creature.'])"
$ws.Range("A5").Value = "('Searing Blaze', ['{R}{R}', 'Instant', 'Searing Blaze deals 1 damage to target player or planeswalker and 1 damage to target creature that player or that planeswalker’s controller controls.', 'Landfall — If you had a land enter the battlefield under your control this turn, Searing Blaze deals 3 damage to that player or planeswalker and 3 damage to that creature instead.'])"
$ws.Range("A6").Value = "('Treasure Hunt', ['{1}{U}', 'Sorcery', 'Reveal cards from the top of your library until you reveal a nonland card, then put all cards revealed this way into your hand.'])"
